$d = $word.ActiveDocument

# Position at the very end of the document content (after "Updated Architecture Design doc.")
$rng = $d.Content
$rng.Collapse(0)

# --- 1. blank centered paragraph ---
$rng.InsertAfter("`r")
$rng.Collapse(0)

# --- 2. "November 17th:" heading paragraph (bold, centered, "th" superscript) ---
$rng.InsertAfter("`r")
$rng.Collapse(0)
$headStart = $rng.Start
$rng.InsertAfter("November 17th:")
$rng.Font.Bold = $true
$rng.Font.BoldBi = $true
$rng.Font.NameAscii = "Times New Roman"
$rng.Font.NameOther = "Times New Roman"
# Make just the "th" superscript, leaving the rest of the run-set untouched.
$thStart = $headStart + ("November 17").Length
$thEnd = $thStart + ("th").Length
$thRange = $d.Range($thStart, $thEnd)
$thRange.Font.Superscript = $true
$rng.Collapse(0)

# --- 3. blank centered paragraph ---
$rng.InsertAfter("`r")
$rng.Collapse(0)

# --- 4. "Brainstormed implementation for friends list in personal profile." ---
$rng.InsertAfter("`r")
$rng.Collapse(0)
$rng.InsertAfter("Brainstormed implementation for friends list in personal profile.")
$rng.Font.NameAscii = "Times New Roman"
$rng.Font.NameOther = "Times New Roman"
$rng.Collapse(0)

# --- 5. blank centered paragraph ---
$rng.InsertAfter("`r")
$rng.Collapse(0)

# --- 6. "Research on unit testing (pt. 2, ft. Jest)." ---
$rng.InsertAfter("`r")
$rng.Collapse(0)
$rng.InsertAfter("Research on unit testing (pt. 2, ft. Jest).")
$rng.Font.NameAscii = "Times New Roman"
$rng.Font.NameOther = "Times New Roman"
$rng.Collapse(0)
